$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

$wsExpo.Range("F7").Value = 4402
$wsExpo.Range("F8").Value = 2627
$wsExpo.Range("F10").Value = 2549
$wsExpo.Range("F15").Value = 672
$wsExpo.Range("F16").Value = 127
$wsExpo.Range("F17").Value = 122
$wsExpo.Range("F18").Value = 335
$wsExpo.Range("F23").Value = 484

# Event cancelled: append marker to the title and mark price as unavailable.
$wsExpo.Range("C25").Value = "杭州·【海潮的回响Echo of The Tide】 | 刀客塔们的大群融入派对·明日方舟SPECIAL ONLY（取消）"
$wsExpo.Range("G25").Value = "不可售"

$wsExpo.Range("F26").Value = 568
$wsExpo.Range("F28").Value = 113
$wsExpo.Range("F29").Value = 420
$wsExpo.Range("F31").Value = 1625
$wsExpo.Range("F32").Value = 1089
$wsExpo.Range("F33").Value = 160
$wsExpo.Range("F35").Value = 1180
$wsExpo.Range("F36").Value = 2087
$wsExpo.Range("F37").Value = 284
$wsExpo.Range("F40").Value = 88
$wsExpo.Range("F43").Value = 672
$wsExpo.Range("F44").Value = 1344
$wsExpo.Range("F45").Value = 119
$wsExpo.Range("F47").Value = 445

# ---------------------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F5").Value = 72

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types)
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F6").Value = 4402
$wsAll.Range("F7").Value = 2627
$wsAll.Range("F8").Value = 2549
$wsAll.Range("F12").Value = 672
$wsAll.Range("F13").Value = 127
$wsAll.Range("F14").Value = 122
$wsAll.Range("F15").Value = 335

# A new event ("EVA ONLY漫展") is inserted at row 19, pushing the two rows
# below it down by one; the former row 21 ("Echo of The Tide") drops out.
$wsAll.Range("C19").Value = "杭州·EVA ONLY漫展"
$wsAll.Range("D19").Value = "文三路199号创业大厦众创空间2层 杭州趣链科技有限公司"
$wsAll.Range("E19").Value = "2024.04.20 10:00-04.20 17:00"
$wsAll.Range("F19").Value = 34
$wsAll.Range("G19").Value = 88
$wsAll.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=82988"
$wsAll.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202403/F9yBq4Qo1710756247458.jpeg"

$wsAll.Range("C20").Value = "杭州·SK怀旧展&偶像专场"
$wsAll.Range("D20").Value = "沈半路171号 T-Car杭州汽车文化主题公园"
$wsAll.Range("E20").Value = "2024.04.20 09:00-04.20 22:00"
$wsAll.Range("F20").Value = 484
$wsAll.Range("G20").Value = 60
$wsAll.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=81764"
$wsAll.Range("I20").Value = "//i1.hdslb.com/bfs/openplatform/202402/mtdbSuTZ1707119415384.jpeg"

$wsAll.Range("C21").Value = "杭州·m字刘海少年和粉毛少女only"
$wsAll.Range("D21").Value = "康候圣街99号 顺丰创新中心"
$wsAll.Range("E21").Value = "2024.04.20 09:00-04.20 17:00"
$wsAll.Range("F21").Value = 28
$wsAll.Range("G21").Value = 68
$wsAll.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=82831"
$wsAll.Range("I21").Value = "//i2.hdslb.com/bfs/openplatform/202403/bVvk6Eky1710383662942.jpeg"

$wsAll.Range("F22").Value = 568
$wsAll.Range("F24").Value = 113
$wsAll.Range("F25").Value = 72
$wsAll.Range("F28").Value = 420
$wsAll.Range("F29").Value = 1625
$wsAll.Range("F30").Value = 1089
$wsAll.Range("F31").Value = 160
$wsAll.Range("F34").Value = 1180
$wsAll.Range("F35").Value = 2087
$wsAll.Range("F36").Value = 284
$wsAll.Range("F41").Value = 88
$wsAll.Range("F43").Value = 672
$wsAll.Range("F44").Value = 1344
$wsAll.Range("F46").Value = 120
$wsAll.Range("F47").Value = 445
